$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.277.92"
$ws.Range("E2").Value = "  +15.25%  "
$ws.Range("D3").Value = "1.683.68"
$ws.Range("E3").Value = "  +9.79%  "
$ws.Range("E4").Value = "  -1.18%  "
$ws.Range("D5").Value = "'307.52"
$ws.Range("E5").Value = "  +9.23%  "
$ws.Range("D6").Value = "'0.9974"
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("D7").Value = "'0.3724"
$ws.Range("E7").Value = "  +2.99%  "
$ws.Range("D8").Value = "'0.3432"
$ws.Range("E8").Value = "  +8.45%  "
$ws.Range("D9").Value = "'48.08"
$ws.Range("E9").Value = "  +18.36%  "
$ws.Range("D10").Value = "'1.184"
$ws.Range("E10").Value = "  +8.57%  "
$ws.Range("D11").Value = "'0.07288"
$ws.Range("E11").Value = "  +7.30%  "
$ws.Range("D12").Value = "'0.9987"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "'20.56"
$ws.Range("E13").Value = "  +10.32%  "
$ws.Range("D14").Value = "'6.089"
$ws.Range("E14").Value = "  +7.81%  "
$ws.Range("D15").Value = "'6.752"
$ws.Range("E15").Value = "  +6.67%  "
$ws.Range("D16").Value = "1.678.28"
$ws.Range("E16").Value = "  +9.81%  "
$ws.Range("D17").Value = "'0.00001106"
$ws.Range("E17").Value = "  +6.35%  "
$ws.Range("D18").Value = "'0.9974"
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("D19").Value = "'0.06711"
$ws.Range("E19").Value = "  +10.48%  "
$ws.Range("D20").Value = "'81.46"
$ws.Range("E20").Value = "  +12.99%  "
$ws.Range("E21").Value = "  +10.05%  "
$ws.Range("D22").Value = "'6.123"
$ws.Range("E22").Value = "  +7.59%  "
$ws.Range("D23").Value = "'12.04"
$ws.Range("E23").Value = "  +6.24%  "
$ws.Range("D24").Value = "24.232.74"
$ws.Range("E24").Value = "  +14.47%  "
$ws.Range("E25").Value = "  +3.06%  "
$ws.Range("D26").Value = "'3.359"
$ws.Range("E26").Value = "  -9.05%  "
$ws.Range("D27").Value = "'2.671"
$ws.Range("E27").Value = "  +21.17%  "
$ws.Range("D28").Value = "'152.32"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("D29").Value = "'19.47"
$ws.Range("E29").Value = "  +10.66%  "
$ws.Range("D30").Value = "1.860.11"
$ws.Range("E30").Value = "  +9.60%  "
$ws.Range("D31").Value = "'126.64"
$ws.Range("E31").Value = "  +6.83%  "
$ws.Range("D32").Value = "'6.414"
$ws.Range("E32").Value = "  +24.39%  "
$ws.Range("D33").Value = "'4.033"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("D34").Value = "'0.9865"
$ws.Range("E34").Value = "  +16.43%  "
$ws.Range("D35").Value = "'1.745"
$ws.Range("E35").Value = "  +15.79%  "
$ws.Range("E36").Value = "  +5.84%  "
$ws.Range("D37").Value = "'12.45"
$ws.Range("E37").Value = "  +17.65%  "
$ws.Range("E38").Value = "  +9.27%  "
$ws.Range("E39").Value = "  +9.41%  "
$ws.Range("D40").Value = "'8.838"
$ws.Range("E40").Value = "  +14.94%  "
$ws.Range("E41").Value = "  +7.44%  "
$ws.Range("D42").Value = "'0.02337"
$ws.Range("E42").Value = "  +11.48%  "
$ws.Range("D43").Value = "'0.2107"
$ws.Range("E43").Value = "  +10.49%  "
$ws.Range("D44").Value = "'0.6140"
$ws.Range("E44").Value = "  +13.28%  "
$ws.Range("D45").Value = "'0.9965"
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("D46").Value = "'3.798"
$ws.Range("E46").Value = "  +6.55%  "
$ws.Range("D47").Value = "'13.17"
$ws.Range("E47").Value = "  +5.47%  "
$ws.Range("D48").Value = "'0.5945"
$ws.Range("E48").Value = "  +9.75%  "
$ws.Range("D49").Value = "'128.20"
$ws.Range("E49").Value = "  +5.64%  "
$ws.Range("D50").Value = "'2.017"
$ws.Range("E50").Value = "  +8.27%  "
$ws.Range("D51").Value = "'0.07155"
$ws.Range("E51").Value = "  +9.14%  "
